# Update "想去人数" (column F) values per the source diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 3032
$ws.Cells.Item(4, 6).Value = 45
$ws.Cells.Item(5, 6).Value = 32
$ws.Cells.Item(7, 6).Value = 163
$ws.Cells.Item(8, 6).Value = 14492
$ws.Cells.Item(9, 6).Value = 160
$ws.Cells.Item(10, 6).Value = 121
$ws.Cells.Item(11, 6).Value = 5795
$ws.Cells.Item(12, 6).Value = 589
$ws.Cells.Item(13, 6).Value = 72
$ws.Cells.Item(14, 6).Value = 43
$ws.Cells.Item(15, 6).Value = 62
$ws.Cells.Item(16, 6).Value = 1237
$ws.Cells.Item(17, 6).Value = 13
$ws.Cells.Item(19, 6).Value = 182
$ws.Cells.Item(20, 6).Value = 793
$ws.Cells.Item(21, 6).Value = 2935
$ws.Cells.Item(22, 6).Value = 54
$ws.Cells.Item(23, 6).Value = 10579
$ws.Cells.Item(24, 6).Value = 1200
$ws.Cells.Item(25, 6).Value = 59
$ws.Cells.Item(26, 6).Value = 82
$ws.Cells.Item(27, 6).Value = 3735
$ws.Cells.Item(28, 6).Value = 244
$ws.Cells.Item(29, 6).Value = 67

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 3032
$ws.Cells.Item(5, 6).Value = 45
$ws.Cells.Item(6, 6).Value = 32
$ws.Cells.Item(8, 6).Value = 165
$ws.Cells.Item(9, 6).Value = 14492
$ws.Cells.Item(10, 6).Value = 160
$ws.Cells.Item(11, 6).Value = 121
$ws.Cells.Item(12, 6).Value = 5795
$ws.Cells.Item(13, 6).Value = 589
$ws.Cells.Item(14, 6).Value = 72
$ws.Cells.Item(15, 6).Value = 43
$ws.Cells.Item(16, 6).Value = 62
$ws.Cells.Item(17, 6).Value = 1237
$ws.Cells.Item(18, 6).Value = 13
$ws.Cells.Item(20, 6).Value = 182
$ws.Cells.Item(21, 6).Value = 793
$ws.Cells.Item(22, 6).Value = 2935
$ws.Cells.Item(23, 6).Value = 54
$ws.Cells.Item(25, 6).Value = 10579
$ws.Cells.Item(26, 6).Value = 1200
$ws.Cells.Item(27, 6).Value = 59
$ws.Cells.Item(28, 6).Value = 82
$ws.Cells.Item(29, 6).Value = 3735
$ws.Cells.Item(30, 6).Value = 244
$ws.Cells.Item(31, 6).Value = 67
